$wb = $excel.ActiveWorkbook

# Sheet 1 "Đơn sale chính" already exists.
$ws1 = $wb.Worksheets.Item(1)

# Add "Đơn thu nợ" right after sheet 1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Đơn thu nợ"

# Add "Lương" right after "Đơn thu nợ".
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Lương"

# ---- Sheet "Đơn sale chính" : header row + 1 order row + "Tổng" totals row ----
$ws1.Cells.Item(1,1).Value = "Tiền tố"
$ws1.Cells.Item(1,2).Value = "Mã dịch vụ"
$ws1.Cells.Item(1,3).Value = "Ngày thực hiện"
$ws1.Cells.Item(1,4).Value = "Cơ sở"
$ws1.Cells.Item(1,5).Value = "Khách hàng"
$ws1.Cells.Item(1,6).Value = "Nguồn khách"
$ws1.Cells.Item(1,7).Value = "Tên dịch vụ"
$ws1.Cells.Item(1,8).Value = "Sale chính"
$ws1.Cells.Item(1,9).Value = "Đơn giá gốc"
$ws1.Cells.Item(1,10).Value = "Sale phụ"
$ws1.Cells.Item(1,11).Value = "Upsale"
$ws1.Cells.Item(1,12).Value = "Đơn giá"
$ws1.Cells.Item(1,13).Value = "Thanh toán lần đầu"
$ws1.Cells.Item(1,14).Value = "Trả sau"
$ws1.Cells.Item(1,15).Value = "Đã thanh toán"
$ws1.Cells.Item(1,16).Value = "Dư nợ"
$ws1.Cells.Item(1,17).Value = "Bác sĩ 1"
$ws1.Cells.Item(1,18).Value = "Bác sĩ 2"
$ws1.Cells.Item(1,19).Value = "Phụ phẫu 1"
$ws1.Cells.Item(1,20).Value = "Phụ phẫu 2"
$ws1.Cells.Item(1,21).Value = "Công phụ phẫu 1"
$ws1.Cells.Item(1,22).Value = "Công phụ phẫu 2"
$ws1.Cells.Item(2,1).Value = "HD-LUXURY"
$ws1.Cells.Item(2,2).Value = 526
$ws1.Cells.Item(2,3).NumberFormat = "@"
$ws1.Cells.Item(2,3).Value = "07-06-2024"
$ws1.Cells.Item(2,4).Value = "CẦN THƠ"
$ws1.Cells.Item(2,5).Value = "Hoàng Thị Thu Vân"
$ws1.Cells.Item(2,6).Value = "CTV"
$ws1.Cells.Item(2,7).Value = "Phun mày"
$ws1.Cells.Item(2,8).Value = "CTV Ngoài"
$ws1.Cells.Item(2,9).Value = 500000
$ws1.Cells.Item(2,10).Value = "Đỗ Thị Huyền Trân"
$ws1.Cells.Item(2,11).Value = 1000000
$ws1.Cells.Item(2,12).Value = 1500000
$ws1.Cells.Item(2,13).Value = 1000000
$ws1.Cells.Item(2,14).Value = 0
$ws1.Cells.Item(2,15).Value = 1000000
$ws1.Cells.Item(2,16).Value = 500000
$ws1.Cells.Item(2,17).Value = "Nguyễn Hoàng Yến Quyên"
$ws1.Cells.Item(2,18).Value = 0
$ws1.Cells.Item(2,19).Value = 0
$ws1.Cells.Item(2,20).Value = 0
$ws1.Cells.Item(2,21).Value = 0
$ws1.Cells.Item(2,22).Value = 0
$ws1.Cells.Item(3,1).Value = "Tổng"
$ws1.Cells.Item(3,2).Value = 1
# C3: blank text cell (left empty)
# D3: blank text cell (left empty)
# E3: blank text cell (left empty)
# F3: blank text cell (left empty)
# G3: blank text cell (left empty)
# H3: blank text cell (left empty)
$ws1.Cells.Item(3,9).Value = 500000
# J3: blank text cell (left empty)
$ws1.Cells.Item(3,11).Value = 1000000
$ws1.Cells.Item(3,12).Value = 1500000
$ws1.Cells.Item(3,13).Value = 1000000
$ws1.Cells.Item(3,14).Value = 0
$ws1.Cells.Item(3,15).Value = 1000000
$ws1.Cells.Item(3,16).Value = 500000
# Q3: blank text cell (left empty)
# R3: blank text cell (left empty)
# S3: blank text cell (left empty)
# T3: blank text cell (left empty)
$ws1.Cells.Item(3,21).Value = 0
$ws1.Cells.Item(3,22).Value = 0

# ---- Sheet "Đơn thu nợ" : header row + 1 debt-collection row + "Tổng" totals row ----
$ws2.Cells.Item(1,1).Value = "Tiền tố"
$ws2.Cells.Item(1,2).Value = "Mã đơn thu nợ"
$ws2.Cells.Item(1,3).Value = "Đơn nợ"
$ws2.Cells.Item(1,4).Value = "Cơ sở"
$ws2.Cells.Item(1,5).Value = "Lượng thu"
$ws2.Cells.Item(1,6).Value = "Sale"
$ws2.Cells.Item(1,7).Value = "Ngày thu"
$ws2.Cells.Item(2,1).Value = "TN"
$ws2.Cells.Item(2,2).Value = 138
$ws2.Cells.Item(2,3).Value = "HD-LUXURY-437"
$ws2.Cells.Item(2,4).Value = "CẦN THƠ"
$ws2.Cells.Item(2,5).Value = 500000
$ws2.Cells.Item(2,6).Value = "CTV Ngoài"
$ws2.Cells.Item(2,7).NumberFormat = "@"
$ws2.Cells.Item(2,7).Value = "07-06-2024"
$ws2.Cells.Item(3,1).Value = "Tổng"
$ws2.Cells.Item(3,2).Value = 1
# C3: blank text cell (left empty)
# D3: blank text cell (left empty)
$ws2.Cells.Item(3,5).Value = 500000
# F3: blank text cell (left empty)
# G3: blank text cell (left empty)

# "Lương" is left blank (just like the placeholder sheet was before this edit).

# Keep the originally active sheet selected.
$ws1.Activate()
